# Insert a new row at position 105 (shifts existing rows 105-117 down to 106-118)
# and populate it with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(105).Insert()

$ws.Cells.Item(105, 1).Value2  = 5
$ws.Cells.Item(105, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(105, 3).Value2  = "Maule"
$ws.Cells.Item(105, 4).Value2  = 45166
$ws.Cells.Item(105, 5).Value2  = 7
$ws.Cells.Item(105, 6).Value2  = 100112040
$ws.Cells.Item(105, 7).Value2  = "Cilantro"
$ws.Cells.Item(105, 8).Value2  = "Sin especificar"
$ws.Cells.Item(105, 9).Value2  = "Primera"
$ws.Cells.Item(105, 10).Value2 = 200
$ws.Cells.Item(105, 11).Value2 = 10000
$ws.Cells.Item(105, 12).Value2 = 10000
$ws.Cells.Item(105, 13).Value2 = 10000
$ws.Cells.Item(105, 14).Value2 = "$/caja 36 atados"
$ws.Cells.Item(105, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(105, 16).Value2 = 278
$ws.Cells.Item(105, 17).Value2 = 36
$ws.Cells.Item(105, 18).Value2 = "Hortaliza"
